# Updated with latest SOP
# Row 33 of the "2025-Master SOP" sheet documents the SOP for mapping
# Functional Tissue Unit illustrations to ASCT+B tables. The SOP title was
# renamed and its version bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2025-Master SOP")

$ws.Range("B33").Value = "Authoring Crosswalk Tables Between Functional Tissue Unit Illustrations and ASCT+B Tables"
$ws.Range("E33").Value = "Version v2.1.1"
